$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the header row (B1:K1) with lowercase, accent-free field names.
# The underlying data columns are unchanged - only the header labels differ.
$ws.Range("B1").Value = "tipo"
$ws.Range("C1").Value = "unidade_dose"
$ws.Range("D1").Value = "dose_minima"
$ws.Range("E1").Value = "dose_maxima"
$ws.Range("F1").Value = "concentracao_maxima"
$ws.Range("G1").Value = "diluicao_sugerida"
$ws.Range("H1").Value = "forma_de_administracao"
$ws.Range("I1").Value = "compativeis"
$ws.Range("J1").Value = "incompativeis"
$ws.Range("K1").Value = "observacoes"

# Format a new (empty) cell J16 with left/center + wrap-text alignment.
$ws.Range("J16").WrapText = $true
$ws.Range("J16").HorizontalAlignment = -4131
$ws.Range("J16").VerticalAlignment = -4108

# Move the active selection to J24.
[void]$ws.Range("J24").Select()
